# Charte de projet - add task detail rows (Priorite, Etat, dates, % acheve, Notes)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Evaluation des charges et calendrier cible
$ws.Range("C9").Value = "Élevée"
$ws.Range("D9").Value = "En cours de réalisation"
$ws.Range("E9").Value = 42468
$ws.Range("F9").Value = 42472
$ws.Range("G9").NumberFormat = "0%"
$ws.Range("G9").Value = 0.5
$ws.Range("I9").Value = "Julie"

# Row 10 - Planification initiale
$ws.Range("C10").Value = "Élevée"
$ws.Range("D10").Value = "En cours de réalisation"
$ws.Range("E10").Value = 42468
$ws.Range("F10").Value = 42472
$ws.Range("G10").NumberFormat = "0%"
$ws.Range("G10").Value = 0.5
$ws.Range("I10").Value = "Julie"

# Row 11 - Gestion du reporting
$ws.Range("C11").Value = "Élevée"
$ws.Range("D11").Value = "Terminée"
$ws.Range("E11").Value = 42468
$ws.Range("F11").Value = 42472
$ws.Range("G11").NumberFormat = "0%"
$ws.Range("G11").Value = 0.75
$ws.Range("I11").Value = "Julie"

# Row 12 - Gestion des relations avec les parties prenantes
$ws.Range("C12").Value = "Élevée"
$ws.Range("D12").Value = "Terminée"
$ws.Range("E12").Value = 42468
$ws.Range("F12").Value = 42472
$ws.Range("G12").NumberFormat = "0%"
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = "Yacine"

# Row 13 - Gestion de la documentation
$ws.Range("C13").Value = "Élevée"
$ws.Range("D13").Value = "Terminée"
$ws.Range("E13").Value = 42468
$ws.Range("F13").Value = 42472
$ws.Range("G13").NumberFormat = "0%"
$ws.Range("G13").Value = 0.75
$ws.Range("I13").Value = "Julie"

# Row 14 - Description des livrables (Priorite/Etat already set)
$ws.Range("E14").Value = 42468
$ws.Range("F14").Value = 42472
$ws.Range("G14").NumberFormat = "0%"
$ws.Range("G14").Value = 1
$ws.Range("I14").Value = "Yacine"

# The Priorité / État drop-down validations used to stop at row 8 then resume
# at row 14-15 (B9:B13 were still blank placeholders). Now that rows 9-13 are
# filled in too, extend both validations to cover the full C5:C15 / D5:D15
# range as one contiguous rule.
$ws.Range("C5:C15").Validation.Delete()
$ws.Range("C5:C15").Validation.Add(3, 1, 1, '"Basse, Normale, Élevée"')
$ws.Range("C5:C15").Validation.ErrorTitle = "Whoops"
$ws.Range("C5:C15").Validation.ShowInput = $false
$ws.Range("C5:C15").Validation.ShowError = $false

$ws.Range("D5:D15").Validation.Delete()
$ws.Range("D5:D15").Validation.Add(3, 2, 1, '"Non commencée,En cours de réalisation, Différé, Terminée"')
$ws.Range("D5:D15").Validation.ErrorTitle = "Whoops"
$ws.Range("D5:D15").Validation.ErrorMessage = "For this template to work correctly you need to select a choice from the drop down list. But you can still use what you entered by clicking Yes."
$ws.Range("D5:D15").Validation.ShowInput = $true
$ws.Range("D5:D15").Validation.ShowError = $true

# Restore the selection the author left the workbook in.
$ws.Range("G9").Select() | Out-Null
